# Generate Report for Handoff
# Updates the generated-file GUID (cd112980-... -> 48164048-...), the
# corresponding xliff hash (bf56d7a... -> 333df60...) and the handoff/
# handback timestamps across the Overview/zh-cn/de-de sheets.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------
$ws_overview.Range("A2").Value = "48164048-4251-40a7-b4cd-cf64c0410084.md"

$hlOverview = $ws_overview.Hyperlinks.Item(1)
$hlOverview.TextToDisplay = "e2e\48164048-4251-40a7-b4cd-cf64c0410084.md"

$ws_overview.Range("G2").Value = "2016-09-02 09:07:18"

# --- zh-cn sheet ------------------------------------------------------
$hlZhCn = $ws_zhcn.Hyperlinks.Item(1)
$hlZhCn.TextToDisplay = "48164048-4251-40a7-b4cd-cf64c0410084.md"

$ws_zhcn.Range("G2").Value = "48164048-4251-40a7-b4cd-cf64c0410084.333df60c181d962e4eb8ecb2e3338662a844bcd8.zh-cn.xlf"
$ws_zhcn.Range("H2").Value = "2016-09-02 09:07:14"

# --- de-de sheet ------------------------------------------------------
$hlDeDe = $ws_dede.Hyperlinks.Item(1)
$hlDeDe.TextToDisplay = "48164048-4251-40a7-b4cd-cf64c0410084.md"

$ws_dede.Range("G2").Value = "48164048-4251-40a7-b4cd-cf64c0410084.333df60c181d962e4eb8ecb2e3338662a844bcd8.de-de.xlf"
$ws_dede.Range("H2").Value = "2016-09-02 09:07:18"
